$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated dSF (column F) values per row, as described in the diff.
$updates = @{
    2  = 3
    3  = 2
    4  = -5
    5  = 3
    6  = -6
    7  = 1
    9  = -1
    10 = 6
    11 = 4
    12 = 1
    13 = 2
    14 = -2
    15 = 1
    16 = -5
    19 = 2
    20 = 5
    21 = 1
    22 = 0
    23 = 2
    24 = 3
    25 = 1
    26 = 7
    27 = -4
    29 = 2
    30 = -5
    32 = 6
    33 = 2
    34 = -4
    35 = -1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
